$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N17").Value = -6002.5002
$ws.Range("L17").Value = 5666.5002
$ws.Range("J17").Value = 1888.8334
$ws.Range("H17").Value = 1891.625
$ws.Range("K31").Value = 675.85716
$ws.Range("I31").Value = 225.28572
$ws.Range("H31").Value = 225.28572
$ws.Range("M31").Value = -445.85716
$ws.Range("L100").Value = 1551.6666
$ws.Range("N100").Value = -2633.6666
$ws.Range("J100").Value = 1551.6666
$ws.Range("H100").Value = 1531.4
$ws.Range("K137").Value = 21432082.5
$ws.Range("I137").Value = 7144027.5
$ws.Range("H137").Value = 1929780
$ws.Range("M137").Value = -21429532.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K45").Value = 85283.60000000001
$ws.Range("M45").Value = -84906.60000000001
$ws.Range("H45").Value = 85283.60000000001
$ws.Range("I45").Value = 85283.60000000001
$ws.Range("N96").Value = -50492
$ws.Range("L96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("H96").Value = 45000
$ws.Range("J122").Value = 0
$ws.Range("H122").Value = 2685.9546
$ws.Range("I122").Value = 2685.9546
$ws.Range("L122").Value = 0
$ws.Range("K122").Value = 8057.8638
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -5607.8638

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3236.1428
$ws.Range("I99").Value = 2610.7646
$ws.Range("K99").Value = 2610.7646
$ws.Range("M99").Value = -1112.7646
$ws.Range("I105").Value = 668281.0600000001
$ws.Range("K105").Value = 668281.0600000001
$ws.Range("M105").Value = -666534.0600000001
$ws.Range("H105").Value = 11306574
$ws.Range("I134").Value = 2658.7585
$ws.Range("H134").Value = 2911.6667
$ws.Range("M134").Value = -5441.2755
$ws.Range("K134").Value = 7976.2755

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("K31").Value = 2357.138
$ws.Range("I31").Value = 2357.138
$ws.Range("H31").Value = 3454.175
$ws.Range("M31").Value = -2062.138
$ws.Range("M34").Value = -2155.138
$ws.Range("H34").Value = 3454.175
$ws.Range("I34").Value = 2357.138
$ws.Range("K34").Value = 2357.138
$ws.Range("N58").Value = -4042
$ws.Range("I58").Value = 2621.8
$ws.Range("J58").Value = 3636
$ws.Range("H58").Value = 3128.9
$ws.Range("L58").Value = 3636
$ws.Range("M58").Value = -2418.8
$ws.Range("K58").Value = 2621.8
$ws.Range("H103").Value = 4903
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("I132").Value = 3529.0557
$ws.Range("L132").Value = 14772.273
$ws.Range("N132").Value = -19832.273
$ws.Range("M132").Value = -8057.167099999999
$ws.Range("H132").Value = 4058.2068
$ws.Range("K132").Value = 10587.1671
$ws.Range("J132").Value = 4924.091
$ws.Range("I136").Value = 2621.8
$ws.Range("K136").Value = 7865.400000000001
$ws.Range("L136").Value = 10908
$ws.Range("H136").Value = 3128.9
$ws.Range("J136").Value = 3636
$ws.Range("N136").Value = -16008
$ws.Range("M136").Value = -5315.400000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N131").Value = -16883.3334
$ws.Range("H131").Value = 9344
$ws.Range("J131").Value = 2267.7778
$ws.Range("L131").Value = 6803.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M80").Value = -111112892
$ws.Range("I80").Value = 111113890
$ws.Range("H80").Value = 71431990
$ws.Range("K80").Value = 111113890
$ws.Range("H83").Value = 71431990
$ws.Range("M83").Value = -555564458
$ws.Range("K83").Value = 555569450
$ws.Range("I83").Value = 111113890
$ws.Range("J102").Value = 2004.3334
$ws.Range("N102").Value = -5248.3334
$ws.Range("I102").Value = 2150
$ws.Range("H102").Value = 2062.6
$ws.Range("L102").Value = 2004.3334
$ws.Range("M102").Value = -528
$ws.Range("K102").Value = 2150
$ws.Range("H113").Value = 1762.3334
$ws.Range("I113").Value = 1521.3846
$ws.Range("K113").Value = 1521.3846
$ws.Range("M113").Value = 648.6153999999999
$ws.Range("H122").Value = 3683.8438
$ws.Range("I122").Value = 2714.261
$ws.Range("K122").Value = 8142.782999999999
$ws.Range("M122").Value = -5692.782999999999
$ws.Range("L133").Value = 160666.33
$ws.Range("N133").Value = -170786.33
$ws.Range("H133").Value = 160666.33
$ws.Range("J133").Value = 160666.33

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("N13").Value = -2780
$ws.Range("J13").Value = 2500
$ws.Range("H16").Value = 863.38464
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340
$ws.Range("M46").Value = -2444.6667
$ws.Range("N46").Value = -2575
$ws.Range("I46").Value = 2632.6667
$ws.Range("H46").Value = 2524.25
$ws.Range("J46").Value = 2199
$ws.Range("L46").Value = 2199
$ws.Range("K46").Value = 2632.6667
$ws.Range("I82").Value = 1341.9117
$ws.Range("K82").Value = 1341.9117
$ws.Range("M82").Value = -980.9117000000001
$ws.Range("H82").Value = 1379.0541
$ws.Range("I85").Value = 1341.9117
$ws.Range("H85").Value = 1379.0541
$ws.Range("M85").Value = -93.91170000000011
$ws.Range("K85").Value = 1341.9117
$ws.Range("H122").Value = 3705.9524
$ws.Range("I122").Value = 3438.3684
$ws.Range("K122").Value = 10315.1052
$ws.Range("M122").Value = -7865.1052
$ws.Range("I136").Value = 2284.2812
$ws.Range("K136").Value = 6852.8436
$ws.Range("H136").Value = 2483.558
$ws.Range("M136").Value = -4302.8436

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1429656.9
$ws.Range("M2").Value = -1429544.9
$ws.Range("K2").Value = 1429656.9
$ws.Range("I2").Value = 1429656.9
$ws.Range("N46").Value = -99959
$ws.Range("H46").Value = 99497
$ws.Range("J46").Value = 99497
$ws.Range("L46").Value = 99497
$ws.Range("N69").ClearContents()
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H122").Value = 20835792
$ws.Range("I122").Value = 2449.7
$ws.Range("K122").Value = 7349.099999999999
$ws.Range("M122").Value = -4899.099999999999
$ws.Range("I132").Value = 6898.4
$ws.Range("L132").Value = 12050.3334
$ws.Range("N132").Value = -17110.3334
$ws.Range("M132").Value = -18165.2
$ws.Range("H132").Value = 5045.9287
$ws.Range("K132").Value = 20695.2
$ws.Range("J132").Value = 4016.7778
$ws.Range("H134").Value = 99497
$ws.Range("J134").Value = 99497
$ws.Range("N134").Value = -303561
$ws.Range("L134").Value = 298491
$ws.Range("I136").Value = 7200.303
$ws.Range("K136").Value = 21600.909
$ws.Range("H136").Value = 7377.9546
$ws.Range("M136").Value = -19050.909
